# Mark exercises 12-16 (rows 12-16, column C) as done (TRUE).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($r in 12..16) {
    $ws.Cells.Item($r, 3).Value = $true
}

# Update the active selection to match the author's final cursor position.
$ws.Range("B11").Select()
